$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77 (pushes old rows 77.. down by one,
# carrying formatting down from the row above as Excel normally does).
$ws.Rows("77").Insert()

# Populate the newly inserted row 77 with the new data record.
$ws.Cells.Item(77, 1).Value = 10
$ws.Cells.Item(77, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(77, 3).Value = "La Araucanía"
$ws.Cells.Item(77, 4).Value = 45167
$ws.Cells.Item(77, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(77, 5).Value = 9
$ws.Cells.Item(77, 6).Value = 100112010
$ws.Cells.Item(77, 7).Value = "Achicoria"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 80
$ws.Cells.Item(77, 11).Value = 10000
$ws.Cells.Item(77, 12).Value = 10000
$ws.Cells.Item(77, 13).Value = 10000
$ws.Cells.Item(77, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(77, 15).Value = "Región Metropolitana"
$ws.Cells.Item(77, 16).Value = 556
$ws.Cells.Item(77, 17).Value = 18
$ws.Cells.Item(77, 18).Value = "Hortaliza"
